# Update cryptos list prices and volume percentages per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.894.78'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.636.40'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.84'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.43'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.78%  '
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.867.20'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.643.13'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.03'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.562'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.888.46'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '229.07'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('E19').Value = '  +2.51%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.05'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.29'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.70%  '
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.55'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('E32').Value = '  +1.27%  '
$ws.Range('E33').Value = '  +1.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.396.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.73%  '
$ws.Range('E35').Value = '  +2.09%  '
$ws.Range('E36').Value = '  +1.92%  '
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.561'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.850'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.55%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('E43').Value = '  +1.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '65.95'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.45'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.774.93'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E47').Value = '  -2.50%  '
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.103'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.74%  '
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.62'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.73%  '
